$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price cells keep their literal text formatting
# (values like "1.00" / "0.608" would otherwise be auto-parsed as numbers by Excel)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.191.53'
$ws.Range("E2").Value = '  +0.73%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.440.94'
$ws.Range("E3").Value = '  +1.49%  '

# Row 4
$ws.Range("E4").Value = '  +0.36%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '406.22'
$ws.Range("E5").Value = '  -2.96%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.02'
$ws.Range("E6").Value = '  +11.30%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.608'
$ws.Range("E7").Value = '  +1.77%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.438.35'
$ws.Range("E8").Value = '  +1.55%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.04%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.680'
$ws.Range("E10").Value = '  +4.73%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.128'
$ws.Range("E11").Value = '  +22.31%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.49'
$ws.Range("E12").Value = '  +3.93%  '

# Row 13
$ws.Range("E13").Value = '  -1.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.984.63'
$ws.Range("E14").Value = '  +1.75%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.65'
$ws.Range("E15").Value = '  +1.24%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.84'
$ws.Range("E16").Value = '  -1.18%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.451.73'
$ws.Range("E17").Value = '  +2.63%  '

# Row 18
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.82'
$ws.Range("E18").Value = '  +9.04%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '62.388.66'
$ws.Range("E19").Value = '  +1.56%  '

# Row 20
$ws.Range("E20").Value = '  -2.05%  '

# Row 21
$ws.Range("E21").Value = '  +19.38%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.27'
$ws.Range("E22").Value = '  -4.48%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '83.10'
$ws.Range("E23").Value = '  +9.17%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.04'
$ws.Range("E24").Value = '  -1.27%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '309.27'
$ws.Range("E25").Value = '  +0.94%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.17'
$ws.Range("E26").Value = '  -2.03%  '

# Row 27
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.54'
$ws.Range("E27").Value = '  +6.74%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '29.83'
$ws.Range("E28").Value = '  +2.66%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.41'
$ws.Range("E29").Value = '  -1.87%  '

# Row 30
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.175'
$ws.Range("E30").Value = '  -2.69%  '

# Row 31
$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.48'
$ws.Range("E31").Value = '  -3.70%  '

# Row 32
$ws.Range("E32").Value = '  +0.45%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '43.25'
$ws.Range("E33").Value = '  +8.58%  '

# Row 34
$ws.Range("E34").Value = '  +0.79%  '

# Row 35
$ws.Range("E35").Value = '  -1.37%  '

# Row 36
$ws.Range("E36").Value = '  +0.03%  '

# Row 37
$ws.Range("E37").Value = '  -5.85%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.58'
$ws.Range("E38").Value = '  -0.18%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.996'
$ws.Range("E39").Value = '  -0.10%  '

# Row 40
$ws.Range("E40").Value = '  +0.63%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.97'
$ws.Range("E41").Value = '  -5.65%  '

# Row 42
$ws.Range("E42").Value = '  +1.64%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.62'
$ws.Range("E43").Value = '  +0.10%  '

# Row 44
$ws.Range("E44").Value = '  +1.23%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.15'
$ws.Range("E45").Value = '  +0.22%  '

# Row 46
$ws.Range("E46").Value = '  -3.60%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.90'
$ws.Range("E47").Value = '  -3.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.25'
$ws.Range("E48").Value = '  -2.27%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.81'
$ws.Range("E49").Value = '  -4.23%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.790.23'
$ws.Range("E50").Value = '  +2.34%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.159.62'
$ws.Range("E51").Value = '  -0.88%  '
